$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Saint Lauren row (row 3) link/copy text
$ws.Range("B3").Value = "Sain Lauren Link"
$ws.Range("C3").Value = "Some copy about Saint Lauren"

# Update Burberry row (row 4)
$ws.Range("B4").Value = "Burberry Link"
$ws.Range("C4").Value = "Some copy about Burberry"

# Update Channel row (row 5)
$ws.Range("B5").Value = "Channel Link"
$ws.Range("C5").Value = "Some copy about Channel"

# Update Louis V row (row 6)
$ws.Range("B6").Value = "Louis V Link"
$ws.Range("C6").Value = "Some copy about Louis V"

# Update MAC row (row 7)
$ws.Range("B7").Value = "MAC Link"
$ws.Range("C7").Value = "Some copy about MAC"

# New row 8: Calvin Klein
$ws.Range("A8").Value = "Calvin Klein"
$ws.Range("B8").Value = "Calvin Klein Link"
$ws.Range("C8").Value = "Some copy about Calvin Klein"

# New row 9: Gucci
$ws.Range("A9").Value = "Gucci"
$ws.Range("B9").Value = "Gucci Link"
$ws.Range("C9").Value = "Some copy about Gucci"

# Copy the style from row 7 (A7:G7) onto rows 8 and 9 so new cells match formatting
# (xlPasteFormats = -4122); values were already set above and are unaffected by a formats-only paste
$ws.Range("A7:G7").Copy()
$ws.Range("A8:G8").PasteSpecial(-4122) | Out-Null
$ws.Range("A9:G9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update row heights
$ws.Rows.Item(3).RowHeight = 32.25
$ws.Rows.Item(4).RowHeight = 32.05
$ws.Rows.Item(5).RowHeight = 32.05
$ws.Rows.Item(6).RowHeight = 32.05
$ws.Rows.Item(7).RowHeight = 32.05
$ws.Rows.Item(8).RowHeight = 32.05
$ws.Rows.Item(9).RowHeight = 32.05
